# Citations slide and citation
#
# Slide 2 (the blank "Title and Content" slide at the end of the deck)
# gets a title of "Citations" and the content placeholder gets a single
# hyperlinked citation line, followed by an empty trailing paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Title placeholder: "Citations" -----------------------------------
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Citations"

# --- Content placeholder: hyperlinked citation + trailing blank line --
$bodyShape = $s.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange

$url = "www.nps.gov/gett/learn/historyculture/civil-war-timeline.htm"

# First paragraph: the citation text.
$bodyTr.Text = $url

# Second (empty) paragraph after it.
[void]$bodyTr.InsertAfter("`r")

# Hyperlink only the first paragraph's run, not the trailing blank one.
$citationRange = $bodyTr.Characters(1, $url.Length)
$citationRange.ActionSettings(1).Hyperlink.Address = "http://www.nps.gov/gett/learn/historyculture/civil-war-timeline.htm"
